$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 22.36000000000006
$ws.Range("H2").Value = [double]"1.468063503636571e-16"
$ws.Range("K2").Value = 59.13781807298047
$ws.Range("L2").Value = "[53.32382998192145, 64.9518061640395]"
$ws.Range("O2").Value = 1.553500271144503
$ws.Range("P2").Value = "[1.4528686746331179, 1.654131867655888]"
$ws.Range("S2").Value = 53.69995897991726
$ws.Range("T2").Value = "[49.69833307010346, 57.70158488973105]"
$ws.Range("W2").Value = 16.83155155155159
$ws.Range("X2").Value = 16.47343343343347
$ws.Range("Y2").Value = 17.18966966966971

# Row 3
$ws.Range("E3").Value = 24.78000000000043
$ws.Range("H3").Value = [double]"1.468063503636571e-16"
$ws.Range("K3").Value = 58.27344929720401
$ws.Range("L3").Value = "[53.1947896474457, 63.35210894696232]"
$ws.Range("O3").Value = -2.289368820634003
$ws.Range("P3").Value = "[-2.377421467581465, -2.201316173686542]"
$ws.Range("S3").Value = 54.90486521619422
$ws.Range("T3").Value = "[52.29284618700027, 57.51688424538818]"
$ws.Range("W3").Value = 9.028948948949109
$ws.Range("X3").Value = 8.681681681681836
$ws.Range("Y3").Value = 9.376216216216383
